$d = $word.ActiveDocument

# Update the date/day heading.
$d.Content.Find.Execute("2023-08-10 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-08-11 Friday", 2)

# Update the division problems in the table, addressed by (row, column)
# so that identical old/new values across different cells can't collide.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "44÷6=7, 2"
$t.Cell(1,2).Range.Text  = "48÷9=5, 3"
$t.Cell(1,3).Range.Text  = "78÷9=8, 6"
$t.Cell(1,4).Range.Text  = "30÷6=5, 0"
$t.Cell(1,5).Range.Text  = "83÷5=16, 3"

$t.Cell(5,1).Range.Text  = "38÷3=12, 2"
$t.Cell(5,2).Range.Text  = "97÷9=10, 7"
$t.Cell(5,3).Range.Text  = "63÷5=12, 3"
$t.Cell(5,4).Range.Text  = "33÷4=8, 1"
$t.Cell(5,5).Range.Text  = "54÷8=6, 6"

$t.Cell(9,1).Range.Text  = "54÷5=10, 4"
$t.Cell(9,2).Range.Text  = "33÷2=16, 1"
$t.Cell(9,3).Range.Text  = "29÷8=3, 5"
$t.Cell(9,4).Range.Text  = "34÷4=8, 2"
$t.Cell(9,5).Range.Text  = "88÷6=14, 4"

$t.Cell(13,1).Range.Text = "96÷4=24, 0"
$t.Cell(13,2).Range.Text = "96÷8=12, 0"
$t.Cell(13,3).Range.Text = "74÷4=18, 2"
$t.Cell(13,4).Range.Text = "83÷9=9, 2"
$t.Cell(13,5).Range.Text = "46÷6=7, 4"

$t.Cell(17,1).Range.Text = "22÷3=7, 1"
$t.Cell(17,2).Range.Text = "26÷5=5, 1"
$t.Cell(17,3).Range.Text = "94÷3=31, 1"
$t.Cell(17,4).Range.Text = "37÷2=18, 1"
$t.Cell(17,5).Range.Text = "22÷8=2, 6"
